$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.075.95"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "3.160.43"
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'580.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").Value = "'149.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.158.02"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").Value = "'6.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").Value = "'37.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "3.679.92"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "64.970.49"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "3.164.59"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'505.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'0.717"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.40%  "
$ws.Range("D23").Value = "'15.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "'7.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").Value = "'84.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").Value = "'9.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("D31").Value = "'27.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'6.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.07%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").Value = "'1.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "'6.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").Value = "'55.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "'0.0889"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("D38").Value = "'479.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("D39").Value = "'0.0416"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.80%  "
$ws.Range("D40").Value = "'2.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").Value = "'8.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("D42").Value = "2.992.43"
$ws.Range("E43").Value = "  -1.91%  "
$ws.Range("E44").Value = "  -3.25%  "
$ws.Range("D45").Value = "'2.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("D46").Value = "'28.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("D47").Value = "0.0₃0596"
$ws.Range("E47").Value = "  +3.89%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("E51").Value = "  +14.38%  "
